# Weekly fruit/vegetable price update.
#
# A new weekly price record is inserted at row 49 (pushing the previously
# existing rows 49-81 down to 50-82, so the sheet grows from A1:R81 to
# A1:R82). The new row carries the same "Vega Modelo de Temuco" /
# "Bruselas (repollito)" series as its neighbours, but with this week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 49:81 down one row, preserving the existing formatting
# (this is what leaves D49 with the date-formatted style once populated).
$ws.Rows(49).Insert()

# Populate the newly opened row 49 with this week's data point.
$ws.Cells.Item(49, 1).Value  = 10
$ws.Cells.Item(49, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(49, 3).Value  = "La Araucanía"
$ws.Cells.Item(49, 4).Value  = 44719
$ws.Cells.Item(49, 5).Value  = 9
$ws.Cells.Item(49, 6).Value  = 100112035
$ws.Cells.Item(49, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Primera"
$ws.Cells.Item(49, 10).Value = 60
$ws.Cells.Item(49, 11).Value = 28000
$ws.Cells.Item(49, 12).Value = 30000
$ws.Cells.Item(49, 13).Value = 29333
$ws.Cells.Item(49, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 2933
$ws.Cells.Item(49, 17).Value = 10
$ws.Cells.Item(49, 18).Value = "Hortaliza"
